$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.471.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.394.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.51%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.399.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.63%  "

$ws.Range("E11").Value = "  +7.47%  "

$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.979.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.49%  "

$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("E15").Value = "  +7.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.538.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.412.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.44%  "

$ws.Range("E20").Value = "  +6.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "386.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("E24").Value = "  +2.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.180"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.50%  "

$ws.Range("E28").Value = "  +17.95%  "

$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("E30").Value = "  +8.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.07%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.89%  "

$ws.Range("E37").Value = "  +10.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0767"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.926.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.18%  "

$ws.Range("E43").Value = "  +5.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.763"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.440.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "298.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.26%  "

$ws.Range("E51").Value = "  -1.98%  "
